# Updates the Epi Info "Survey" Excel template:
#  - adds two new list-sheets ("Sheet4" = Yes/No/Don't know, "Sheet2" = Foods Eaten options)
#  - adds one row to the DataTypes list ("Time")
#  - adds three new question rows to Sheet1 (Symptoms Time, Foods Eaten, Hospitalization)
#  - tweaks the Title/Description text for the existing "symptoms start date" question

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert the two new lookup-list sheets.
#    Excel's auto-naming picks the lowest unused "SheetN" number; inserting
#    both right after Sheet1 (in this order) reproduces the target tab order
#    Sheet1, Sheet4, Sheet2, Sheet3, DataTypes.
#    NB: re-fetch the new sheets *by name* afterwards instead of keeping the
#    references returned by Add() - those track sheet *position*, which
#    shifts when the second sheet is inserted.
# ---------------------------------------------------------------------------
$sheet1 = $wb.Worksheets.Item("Sheet1")

[void]$wb.Worksheets.Add($null, $sheet1)   # -> auto-named "Sheet2"
[void]$wb.Worksheets.Add($null, $sheet1)   # -> auto-named "Sheet4"

# Re-fetch every sheet reference *by name* now that the sheet collection is
# done changing shape - Add()'s return value (and any reference captured
# beforehand) tracks position, not sheet identity, so it can silently point
# at the wrong tab once more sheets are inserted/reordered.
$sheet1 = $wb.Worksheets.Item("Sheet1")
$sheetYesNo = $wb.Worksheets.Item("Sheet4")
$sheetFoods = $wb.Worksheets.Item("Sheet2")

# "Sheet4" tab: Yes / No / Don't know  (used by the new Hospitalization question)
$sheetYesNo.Range("A1").Value = "Yes"
$sheetYesNo.Range("A2").Value = "No"
$sheetYesNo.Range("A3").Value = "Don't know "

# "Sheet2" tab: list of foods (used by the new Foods Eaten question)
$sheetFoods.Range("A1").Value = "Fresh celery"
$sheetFoods.Range("A2").Value = "Grapes"
$sheetFoods.Range("A3").Value = "Peaches"
$sheetFoods.Range("A4").Value = "Apple juice"
$sheetFoods.Range("A5").Value = "Orange juice"

# ---------------------------------------------------------------------------
# 2. DataTypes sheet: add a "Time" data type under the existing list.
# ---------------------------------------------------------------------------
$dataTypes = $wb.Worksheets.Item("DataTypes")
$dataTypes.Range("A8").Value = "Time"

# ---------------------------------------------------------------------------
# 3. Sheet1: update the symptoms-start-date row and append three new
#    question rows (Symptoms Time, Foods Eaten, Hospitalization).
# ---------------------------------------------------------------------------

# Row 6 ("When did symptoms start?") - title/description become date-specific
# now that a separate time question exists.
$sheet1.Range("B6").Value = "Symptoms Date"
$sheet1.Range("C6").Value = "Please enter symptoms date…"

# Row 7 (new): Symptoms Time
$sheet1.Range("A7").Value = "What time did the symptoms start?"
$sheet1.Range("B7").Value = "Symptoms Time"
$sheet1.Range("C7").Value = "Please enter symptoms time…"
$sheet1.Range("D7").Value = "onse_time"
$sheet1.Range("E7").Value = "Time"
$sheet1.Range("F7").Value = $false

# Row 8 (new): Foods Eaten
$sheet1.Range("A8").Value = "Select eaten foods:"
$sheet1.Range("B8").Value = "Foods Eaten"
$sheet1.Range("C8").Value = "Please select eaten foods…"
$sheet1.Range("D8").Value = "eaten_foods"
$sheet1.Range("E8").Value = "Checkbox"
$sheet1.Range("F8").Value = $false
$sheet1.Range("G8").Value = "Sheet2"

# Row 9 (new): Hospitalization
$sheet1.Range("A9").Value = "Was patient hospitalized?"
$sheet1.Range("B9").Value = "Hospitalization"
$sheet1.Range("D9").Value = "hospitalized"
$sheet1.Range("E9").Value = "Options"
$sheet1.Range("F9").Value = $false
$sheet1.Range("G9").Value = "Sheet4"
